$p = $ppt.ActivePresentation

# Append a new slide (17th) using the "Title and Content" layout (2 = ppLayoutText),
# matching the layout used by the other content slides (slideLayout2.xml).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Pushing Changes"

# Content placeholder (two paragraphs)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Say I created a new file – maxbrown.txt `rI now want that new file stored in my repo"
